$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Value = 44231
$ws.Range("B15").Value = 0.16666666666666666
$ws.Range("C15").Value = "Procurando solução para assossiar Produtos com Grupos"

$ws.Range("A16").Value = 44235
$ws.Range("B16").Value = 0.29166666666666669
$ws.Range("C16").Value = "Procurado e encontrado solução p/associar Produto c/Grupo"

$ws.Range("A17").Value = 44236
$ws.Range("B17").Value = 0.16666666666666666
$ws.Range("C17").Value = "Procurando solução para assossiar Pedidos com Clientes"

$ws.Range("C18").Select()
